# Update the ENVELOPE_ASSEMBLIES sheet: rename the short construction-type
# codes (WALL2, SHADE2, ROOF2, WIN4, LEAK3, CONS3, ...) to their new
# descriptive names (WALL_AS2, SHADING_AS2, ROOF_AS2, WINDOW_AS4,
# TIGHTNESS_AS3, CONSTRUCTION_AS3, ...).
#
# Columns are touched in this exact order (F, E, G, D, C, B) and each
# column top-to-bottom, so that the newly-introduced shared strings are
# interned in the same order they appear in the refreshed workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")

# Column F - Wall construction type
$ws.Cells.Item(2, 6).Value = "WALL_AS2"
$ws.Cells.Item(3, 6).Value = "WALL_AS5"
$ws.Cells.Item(4, 6).Value = "WALL_AS5"
$ws.Cells.Item(5, 6).Value = "WALL_AS5"
$ws.Cells.Item(6, 6).Value = "WALL_AS5"
$ws.Cells.Item(7, 6).Value = "WALL_AS5"

# Column E - Roof construction type
$ws.Cells.Item(2, 5).Value = "ROOF_AS2"
$ws.Cells.Item(3, 5).Value = "ROOF_AS2"
$ws.Cells.Item(4, 5).Value = "ROOF_AS2"
$ws.Cells.Item(5, 5).Value = "ROOF_AS1"
$ws.Cells.Item(6, 5).Value = "ROOF_AS4"
$ws.Cells.Item(7, 5).Value = "ROOF_AS4"

# Column G - Shading type
$ws.Cells.Item(2, 7).Value = "SHADING_AS2"
$ws.Cells.Item(3, 7).Value = "SHADING_AS1"
$ws.Cells.Item(4, 7).Value = "SHADING_AS1"
$ws.Cells.Item(5, 7).Value = "SHADING_AS1"
$ws.Cells.Item(6, 7).Value = "SHADING_AS1"
$ws.Cells.Item(7, 7).Value = "SHADING_AS1"

# Column D - Window type
$ws.Cells.Item(2, 4).Value = "WINDOW_AS1"
$ws.Cells.Item(3, 4).Value = "WINDOW_AS1"
$ws.Cells.Item(4, 4).Value = "WINDOW_AS1"
$ws.Cells.Item(5, 4).Value = "WINDOW_AS2"
$ws.Cells.Item(6, 4).Value = "WINDOW_AS2"
$ws.Cells.Item(7, 4).Value = "WINDOW_AS4"

# Column C - Air tightness type
$ws.Cells.Item(2, 3).Value = "TIGHTNESS_AS3"
$ws.Cells.Item(3, 3).Value = "TIGHTNESS_AS3"
$ws.Cells.Item(4, 3).Value = "TIGHTNESS_AS3"
$ws.Cells.Item(5, 3).Value = "TIGHTNESS_AS3"
$ws.Cells.Item(6, 3).Value = "TIGHTNESS_AS2"
$ws.Cells.Item(7, 3).Value = "TIGHTNESS_AS1"

# Column B - Construction type
$ws.Cells.Item(2, 2).Value = "CONSTRUCTION_AS3"
$ws.Cells.Item(3, 2).Value = "CONSTRUCTION_AS2"
$ws.Cells.Item(4, 2).Value = "CONSTRUCTION_AS2"
$ws.Cells.Item(5, 2).Value = "CONSTRUCTION_AS3"
$ws.Cells.Item(6, 2).Value = "CONSTRUCTION_AS3"
$ws.Cells.Item(7, 2).Value = "CONSTRUCTION_AS3"

# Column widths grew now that the codes are longer -- refresh the
# best-fit widths for the touched columns.
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null
$ws.Columns.Item(7).AutoFit() | Out-Null

# The workbook was last saved with this sheet active (cell N13 selected),
# moving the "last active" tab away from SUPPLY_ASSEMBLIES.
$ws.Activate()
$ws.Range("N13").Select()
